$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 75, shifting existing rows 75-126 down to 76-127
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new record's data
$ws.Range("A75").Value = 3
$ws.Range("B75").Value = "Femacal de La Calera"
$ws.Range("C75").Value = "Coquimbo"
$ws.Range("D75").Value = 44574
$ws.Range("E75").Value = 5
$ws.Range("F75").Value = 100112030
$ws.Range("G75").Value = "Poroto granado"
$ws.Range("H75").Value = "Sin especificar"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 100
$ws.Range("K75").Value = 28000
$ws.Range("L75").Value = 28000
$ws.Range("M75").Value = 28000
$ws.Range("N75").Value = "`$/saco 25 kilos"
$ws.Range("O75").Value = "Provincia de Quillota"
$ws.Range("P75").Value = 1120
$ws.Range("Q75").Value = 25
$ws.Range("R75").Value = "Hortaliza"
